# Updates the cryptos list worksheet in-place with refreshed price/volume
# data (and, for a few coins, swapped ranking positions), per the commit:
# "Updated cryptos list on Sat Apr 15 13:21:07 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write the "Price" column (D) as literal text. Some price strings
# (e.g. "334.47") look like plain numbers to Excel's automatic type
# detection and would otherwise be stored as a numeric value instead of
# text, so force those with a leading apostrophe (exactly what Excel does
# when a user types a number-looking value into a text entry).
function Set-Price($row, $val) {
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Range("D$row").Value = "'" + $val
    } else {
        $ws.Range("D$row").Value = $val
    }
}

function Set-Volume($row, $val) {
    $ws.Range("E$row").Value = $val
}

function Set-Coin($row, $name, $link) {
    $ws.Range("B$row").Value = $name
    $ws.Range("C$row").Value = $link
}

# Row 2: Bitcoin
Set-Price  2 "30.484.52"
Set-Volume 2 "  -0.87%  "

# Row 3: Ethereum
Set-Price  3 "2.109.81"
Set-Volume 3 "  +0.17%  "

# Row 4: TetherUSD
Set-Volume 4 "  -0.10%  "

# Row 5: BNB
Set-Price 5 "334.47"

# Row 6: USDC
Set-Price  6 "1.001"
Set-Volume 6 "  -0.10%  "

# Row 7: XRP
Set-Price  7 "0.5255"
Set-Volume 7 "  -0.90%  "

# Row 8: Cardano
Set-Price  8 "0.4572"
Set-Volume 8 "  +4.66%  "

# Row 9: OKB
Set-Price  9 "53.64"
Set-Volume 9 "  +14.27%  "

# Row 10: Dogecoin
Set-Price  10 "0.08989"
Set-Volume 10 "  +0.46%  "

# Row 11: Polygon
Set-Price  11 "1.181"
Set-Volume 11 "  +1.25%  "

# Row 12: Solana
Set-Price  12 "24.46"
Set-Volume 12 "  -1.28%  "

# Row 13: WrappedEther
Set-Price  13 "2.105.69"
Set-Volume 13 "  +0.11%  "

# Row 14: Polkadot
Set-Price  14 "6.821"
Set-Volume 14 "  +1.38%  "

# Row 15: Chainlink
Set-Price  15 "7.842"
Set-Volume 15 "  +1.06%  "

# Row 16: Litecoin
Set-Price  16 "96.81"
Set-Volume 16 "  +0.20%  "

# Row 17: BinanceUSD
Set-Price  17 "1.004"
Set-Volume 17 "  -0.02%  "

# Row 18: ShibaInu
Set-Volume 18 "  +0.03%  "

# Row 19: TRON
Set-Price  19 "0.06617"
Set-Volume 19 "  -1.05%  "

# Row 20: Avalanche
Set-Price  20 "19.41"
Set-Volume 20 "  +2.11%  "

# Row 21: Dai
Set-Volume 21 "  -0.14%  "

# Row 22: Uniswap
Set-Price  22 "6.306"
Set-Volume 22 "  +0.13%  "

# Row 23: WrappedBTC
Set-Price 23 "30.557.82"

# Row 24: Cosmos
Set-Price  24 "12.36"
Set-Volume 24 "  +0.82%  "

# Row 25: Toncoin
Set-Price  25 "2.346"
Set-Volume 25 "  +2.55%  "

# Row 26: WrappedliquidstakedEther2.0
Set-Price  26 "2.346.21"
Set-Volume 26 "  -0.26%  "

# Row 27: EthereumClassic
Set-Price  27 "22.41"
Set-Volume 27 "  -0.84%  "

# Row 28: LidoDAOToken
Set-Price  28 "2.582"
Set-Volume 28 "  +0.23%  "

# Row 29: Monero
Set-Price  29 "163.46"
Set-Volume 29 "  +0.43%  "

# Row 30: BitcoinCash
Set-Price  30 "132.86"
Set-Volume 30 "  +0.01%  "

# Row 31: ImmutableX
Set-Price  31 "1.197"
Set-Volume 31 "  +0.40%  "

# Row 32 & 33 swap ranking positions: Stellar <-> ARBITRUM
Set-Coin   32 "ARBITRUM" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-Price  32 "1.719"
Set-Volume 32 "  +11.63%  "

Set-Coin   33 "Stellar" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-Price  33 "0.1074"
Set-Volume 33 "  -0.60%  "

# Row 34: Filecoin
Set-Price  34 "6.162"
Set-Volume 34 "  -0.14%  "

# Row 35: HuobiToken
Set-Price  35 "3.926"
Set-Volume 35 "  -2.70%  "

# Row 36: FraxShare
Set-Price  36 "10.44"
Set-Volume 36 "  +9.06%  "

# Row 37: VeChain
Set-Price  37 "0.02581"
Set-Volume 37 "  -0.56%  "

# Row 38 & 39 swap ranking positions: InternetComputer(DFINITY) <-> Hedera
Set-Coin   38 "Hedera" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Price  38 "0.06827"
Set-Volume 38 "  +1.17%  "

Set-Coin   39 "InternetComputer(DFINITY)" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-Price  39 "5.560"
Set-Volume 39 "  +0.50%  "

# Row 40: Aptos
Set-Price  40 "12.78"
Set-Volume 40 "  +1.04%  "

# Row 41: Algorand
Set-Price  41 "0.2298"
Set-Volume 41 "  +1.20%  "

# Row 42: TheSandbox
Set-Price  42 "0.6923"
Set-Volume 42 "  +1.66%  "

# Row 43: TrustWalletToken
Set-Price  43 "1.247"
Set-Volume 43 "  +0.11%  "

# Row 44: NEARProtocol
Set-Price  44 "2.356"
Set-Volume 44 "  +6.32%  "

# Row 45: Frax
Set-Volume 45 "  -0.07%  "

# Row 46 & 47 swap ranking positions: Decentraland <-> EnergySwap
Set-Coin   46 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Price  46 "14.09"
Set-Volume 46 "  +0.14%  "

Set-Coin   47 "Decentraland" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-Price  47 "0.6391"
Set-Volume 47 "  -0.27%  "

# Row 48: PancakeSwap
Set-Price  48 "3.658"
Set-Volume 48 "  -0.11%  "

# Row 49: BabyDogeCoin
Set-Price  49 "0.00000000357"
Set-Volume 49 "  +25.46%  "

# Row 50: EOS
Set-Price  50 "1.250"
Set-Volume 50 "  -0.62%  "

# Row 51: WEMIXTOKEN
Set-Price  51 "1.222"
Set-Volume 51 "  +2.63%  "
